$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '76.022.78'
Set-TextValue 'E2' '  +1.41%  '
Set-TextValue 'D3' '2.916.39'
Set-TextValue 'E3' '  +3.72%  '
Set-TextValue 'E4' '  +0.04%  '
Set-TextValue 'D5' '203.64'
Set-TextValue 'E5' '  +8.92%  '
Set-TextValue 'D6' '596.10'
Set-TextValue 'E6' '  +0.22%  '
Set-TextValue 'E7' '  -0.04%  '
Set-TextValue 'E8' '  +0.40%  '
Set-TextValue 'D9' '0.197'
Set-TextValue 'E9' '  +2.78%  '
Set-TextValue 'D10' '2.912.63'
Set-TextValue 'E10' '  +3.61%  '
Set-TextValue 'D11' '0.439'
Set-TextValue 'E11' '  +17.73%  '
Set-TextValue 'D12' '0.162'
Set-TextValue 'E12' '  +0.71%  '
Set-TextValue 'E13' '  +0.45%  '
Set-TextValue 'D14' '3.452.51'
Set-TextValue 'E14' '  +3.73%  '
Set-TextValue 'D15' '28.35'
Set-TextValue 'E15' '  +5.52%  '
Set-TextValue 'D16' '75.951.87'
Set-TextValue 'E16' '  +1.54%  '
Set-TextValue 'E17' '  +1.09%  '
Set-TextValue 'D18' '2.907.40'
Set-TextValue 'E18' '  +3.43%  '
Set-TextValue 'D19' '13.08'
Set-TextValue 'E19' '  +6.42%  '
Set-TextValue 'D20' '8.88'
Set-TextValue 'E20' '  -1.48%  '
Set-TextValue 'D21' '373.24'
Set-TextValue 'E21' '  -1.22%  '
Set-TextValue 'E22' '  +1.28%  '
Set-TextValue 'D23' '4.31'
Set-TextValue 'E23' '  +5.44%  '
Set-TextValue 'D24' '71.38'
Set-TextValue 'E24' '  +0.51%  '
Set-TextValue 'E25' '  -0.05%  '
Set-TextValue 'E26' '  +3.81%  '
Set-TextValue 'E27' '  +2.33%  '
Set-TextValue 'D28' '9.71'
Set-TextValue 'E28' '  -1.37%  '
Set-TextValue 'E29' '  +3.91%  '
Set-TextValue 'D30' '0.998'
Set-TextValue 'E30' '  -0.16%  '
Set-TextValue 'E31' '  +0.23%  '
Set-TextValue 'D32' '503.14'
Set-TextValue 'E32' '  -3.48%  '
Set-TextValue 'D33' '7.79'
Set-TextValue 'E33' '  +1.64%  '
Set-TextValue 'E34' '  +2.83%  '
Set-TextValue 'E35' '  +0.11%  '
Set-TextValue 'B36' 'Monero'
Set-TextValue 'C36' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D36' '165.21'
Set-TextValue 'E36' '  +1.33%  '
Set-TextValue 'B37' 'EthereumClassic'
Set-TextValue 'C37' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D37' '20.23'
Set-TextValue 'E37' '  +1.39%  '
Set-TextValue 'D38' '0.109'
Set-TextValue 'E38' '  +26.85%  '
Set-TextValue 'D39' '19.60'
Set-TextValue 'E39' '  +1.13%  '
Set-TextValue 'B40' 'PolygonEcosystemToken'
Set-TextValue 'C40' 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 'D40' '0.376'
Set-TextValue 'E40' '  +10.03%  '
Set-TextValue 'B41' 'Kaspa'
Set-TextValue 'C41' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D41' '0.113'
Set-TextValue 'E41' '  -4.98%  '
Set-TextValue 'D42' '182.42'
Set-TextValue 'E42' '  -1.88%  '
Set-TextValue 'E43' '  -0.06%  '
Set-TextValue 'E44' '  +0.07%  '
Set-TextValue 'E45' '  -0.21%  '
Set-TextValue 'D46' '40.15'
Set-TextValue 'E46' '  +0.45%  '
Set-TextValue 'D47' '1.20'
Set-TextValue 'E47' '  -1.78%  '
Set-TextValue 'D48' '2.34'
Set-TextValue 'E48' '  +0.80%  '
Set-TextValue 'D49' '0.574'
Set-TextValue 'E49' '  -0.77%  '
Set-TextValue 'B50' 'Filecoin'
Set-TextValue 'C50' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D50' '3.72'
Set-TextValue 'E50' '  +0.18%  '
Set-TextValue 'B51' 'InjectiveProtocol'
Set-TextValue 'C51' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D51' '22.45'
Set-TextValue 'E51' '  +7.08%  '
